$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of the TC2 step 2.0 row (row 20) and the TC3 step 2.0 row (row 28)
$ws.Range("B20").Value = "Chefe Clica para realizar a liquidação."
$ws.Range("D20").Value = "SYSTEM Apresenta a tela de Registrar Liquidações"

$ws.Range("B28").Value = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
